# Update cryptocurrency price/volume snapshot (GitHub Actions scheduled refresh).
# For D-column cells whose new value would otherwise be auto-parsed by Excel
# as a genuine number (stripping trailing zeros / changing representation),
# we force a text/quote-prefix entry so the stored cell keeps the exact
# decimal-string formatting used throughout this sheet (mirrors the original
# data, which was written as literal text, not numeric, values).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '59.378.43'
$ws.Range('E2').Value = '  -0.08%  '
$ws.Range('D3').Value = '2.516.95'
$ws.Range('E3').Value = '  -0.41%  '
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').Value = "'537.48"
$ws.Range('E5').Value = '  -1.20%  '
$ws.Range('D6').Value = "'139.53"
$ws.Range('E6').Value = '  -4.12%  '
$ws.Range('D7').Value = "'0.999"
$ws.Range('E7').Value = '  +0.35%  '
$ws.Range('E8').Value = '  -1.78%  '
$ws.Range('D9').Value = '2.521.29'
$ws.Range('E9').Value = '  -1.18%  '
$ws.Range('E10').Value = '  +0.52%  '
$ws.Range('E11').Value = '  +1.01%  '
$ws.Range('E12').Value = '  -3.12%  '
$ws.Range('D13').Value = "'0.358"
$ws.Range('E13').Value = '  -1.40%  '
$ws.Range('D14').Value = '2.965.73'
$ws.Range('E14').Value = '  -0.26%  '
$ws.Range('D15').Value = "'23.48"
$ws.Range('E15').Value = '  -1.60%  '
$ws.Range('D16').Value = '59.281.04'
$ws.Range('E16').Value = '  -0.07%  '
$ws.Range('E17').Value = '  -0.44%  '
$ws.Range('D18').Value = '2.517.40'
$ws.Range('E18').Value = '  -0.74%  '
$ws.Range('D19').Value = "'11.11"
$ws.Range('E19').Value = '  -1.02%  '
$ws.Range('D20').Value = "'4.32"
$ws.Range('E20').Value = '  +0.51%  '
$ws.Range('D21').Value = "'325.89"
$ws.Range('E21').Value = '  -0.27%  '
$ws.Range('D22').Value = "'0.999"
$ws.Range('E22').Value = '  +0.08%  '
$ws.Range('D24').Value = "'63.21"
$ws.Range('E24').Value = '  +1.92%  '
$ws.Range('D25').Value = "'0.425"
$ws.Range('E25').Value = '  -2.21%  '
$ws.Range('E26').Value = '  +1.94%  '
$ws.Range('E27').Value = '  +0.79%  '
$ws.Range('D28').Value = "'7.83"
$ws.Range('E28').Value = '  -2.00%  '
$ws.Range('D29').Value = "'6.95"
$ws.Range('E29').Value = '  +3.44%  '
$ws.Range('E30').Value = '  -0.88%  '
$ws.Range('D31').Value = "'1.79"
$ws.Range('E31').Value = '  -2.27%  '
$ws.Range('D32').Value = "'164.82"
$ws.Range('E32').Value = '  +2.80%  '
$ws.Range('E33').Value = '  +0.18%  '
$ws.Range('E34').Value = '  -3.06%  '
$ws.Range('E35').Value = '  -7.00%  '
$ws.Range('D36').Value = "'18.53"
$ws.Range('E36').Value = '  -1.11%  '
$ws.Range('D37').Value = "'4.27"
$ws.Range('E37').Value = '  -3.26%  '
$ws.Range('E38').Value = '  -1.73%  '
$ws.Range('D39').Value = "'36.95"
$ws.Range('E39').Value = '  -0.12%  '
$ws.Range('E40').Value = '  -0.77%  '
$ws.Range('D41').Value = "'0.812"
$ws.Range('E41').Value = '  -2.57%  '
$ws.Range('D42').Value = "'5.24"
$ws.Range('E42').Value = '  -6.82%  '
$ws.Range('D43').Value = "'279.81"
$ws.Range('E43').Value = '  -5.16%  '
$ws.Range('D44').Value = "'0.998"
$ws.Range('E44').Value = '  +0.50%  '
$ws.Range('E45').Value = '  -0.69%  '
$ws.Range('D46').Value = "'10.88"
$ws.Range('E46').Value = '  +0.60%  '
$ws.Range('D47').Value = "'0.0933"
$ws.Range('E47').Value = '  -0.38%  '
$ws.Range('D48').Value = "'123.64"
$ws.Range('E48').Value = '  +0.20%  '
$ws.Range('E49').Value = '  -0.24%  '
$ws.Range('E50').Value = '  -1.68%  '
$ws.Range('D51').Value = "'17.88"
$ws.Range('E51').Value = '  -2.16%  '
